$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9677245616912842
$ws.Range("B1").Value = 1.309038758277893
$ws.Range("C1").Value = 2.202776193618774
$ws.Range("D1").Value = -1
$ws.Range("E1").Value = 1.97918975353241
